$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A13: new rich-text cell (test string exercising bold/italic/underline
# and multi-space / period runs, per "work with <br>" fix) -------------------
$cell = $ws.Range("A13")

# The row previously just held an empty, wrap-text styled cell (s="1").
# The new cell carries its own per-run fonts and no longer needs that style.
$cell.ClearFormats() | Out-Null

$text = "test này là in đậm.    , test.     này là in      nghiệm.   , test này là underline"
$cell.Value = $text

# Run 1 "test này là "          -> left as default (no explicit formatting)
# Run 2 "in đậm.    "           -> bold
$cell.Characters(13, 11).Font.Bold = $true
# Run 3 ", test.     này là "   -> plain run, but with explicit font info
$cell.Characters(24, 19).Font.Size = 12
# Run 4 "in      nghiệm.   , "  -> italic
$cell.Characters(43, 20).Font.Italic = $true
# Run 5 "test này là "          -> plain run, but with explicit font info
$cell.Characters(63, 12).Font.Size = 12
# Run 6 "underline"             -> underline, body font
$cell.Characters(75, 9).Font.Underline = $true
$cell.Characters(75, 9).Font.Name = "Aptos Narrow (Body)"

# --- sheet view: scroll back to top, move selection to B6 -------------------
$ws.Range("B6").Select() | Out-Null

Write-Output "done"
